$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking values
# (e.g. "1.003") are preserved as text, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.178.11"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "1.800.67"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "314.57"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "0.5294"
$ws.Range("E7").Value = "  +3.70%  "
$ws.Range("D8").Value = "0.3831"
$ws.Range("E8").Value = "  -2.41%  "
$ws.Range("D9").Value = "0.07997"
$ws.Range("E9").Value = "  +2.82%  "
$ws.Range("D10").Value = "41.34"
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("D12").Value = "6.332"
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").Value = "1.003"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").Value = "20.60"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").Value = "1.804.05"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("E16").Value = "  -1.87%  "
$ws.Range("D17").Value = "92.73"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "0.00001096"
$ws.Range("E18").Value = "  -3.55%  "
$ws.Range("D19").Value = "0.06609"
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("E21").Value = "  -1.92%  "
$ws.Range("D22").Value = "5.975"
$ws.Range("E22").Value = "  -1.73%  "
$ws.Range("D23").Value = "28.221.39"
$ws.Range("E23").Value = "  -0.74%  "
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D26").Value = "160.27"
$ws.Range("E26").Value = "  +3.37%  "
$ws.Range("E27").Value = "  -2.99%  "
$ws.Range("D28").Value = "2.011.59"
$ws.Range("E28").Value = "  -0.63%  "
$ws.Range("D29").Value = "2.386"
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("D30").Value = "123.25"
$ws.Range("E30").Value = "  -1.65%  "
$ws.Range("D31").Value = "0.1091"
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("E32").Value = "  -3.74%  "
$ws.Range("D33").Value = "3.656"
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("D34").Value = "5.556"
$ws.Range("E34").Value = "  -1.79%  "
$ws.Range("D35").Value = "0.07284"
$ws.Range("E35").Value = "  +3.45%  "
$ws.Range("D36").Value = "12.29"
$ws.Range("E36").Value = "  +10.03%  "
$ws.Range("D37").Value = "8.933"
$ws.Range("E37").Value = "  +2.18%  "
$ws.Range("D38").Value = "0.2167"
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("D39").Value = "0.02316"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "5.080"
$ws.Range("E40").Value = "  -1.99%  "
$ws.Range("D41").Value = "0.6206"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("D42").Value = "1.164"
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("D43").Value = "1.370"
$ws.Range("E43").Value = "  -0.97%  "
$ws.Range("D44").Value = "13.29"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D45").Value = "0.6000"
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("D46").Value = "3.763"
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("E47").Value = "  +2.07%  "
$ws.Range("D48").Value = "1.210"
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("D49").Value = "1.929"
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("D50").Value = "0.06831"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").Value = "73.11"
$ws.Range("E51").Value = "  -1.24%  "
